$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '242.11'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '21.82'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.393'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05724'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '3.434'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '6.312'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8076'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.8568'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1439'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03041'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03122'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09359'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.928'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.001587'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.04810'
$ws.Range("B18").Value = 'One'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0005865'
$ws.Range("E18").Value = '17OneONE'
$ws.Range("B19").Value = 'TigerCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.006405'
$ws.Range("E19").Value = '18TigerCashTCH'
$ws.Range("B20").Value = 'HotbitToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.004091'
$ws.Range("E20").Value = '19HotbitTokenHTB'
$ws.Range("B21").Value = 'BitKan'
$ws.Range("C21").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.001001'
$ws.Range("E21").Value = '20BitKanKAN'
$ws.Range("B22").Value = 'NitroEx'
$ws.Range("C22").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0001504'
$ws.Range("E22").Value = '21NitroExNTX'
$ws.Range("B23").Value = 'LEO'
$ws.Range("C23").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.724'
$ws.Range("E23").Value = '22LEOLEO'
$ws.Range("B24").Value = 'BTSEToken'
$ws.Range("C24").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.176'
$ws.Range("E24").Value = '23BTSETokenBTSE'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1254'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0004010'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006720'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1057'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002564'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.007337'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005632'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000752'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5814'
$ws.Range("E47").Value = '46CoinbaseStockTokenCOINBestin24h'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.1423'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00002105'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.01012'
